$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row at row 5 (shifts the existing rows 5..190 down to 6..191,
# extending formulas / hyperlinks / dimension exactly like a native Excel
# "Insert Row" does).
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Insert()

# ---------------------------------------------------------------------------
# Populate the freshly inserted row 5 with the new
# general.hotelLocationFileName configuration entry.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "*"
$ws.Range("B5").Value = "general"
$ws.Range("C5").Value = "hotelLocationFileName"
$ws.Range("D5").Formula = '="@@."&A5&"."&B5&"."&C5&"@@"'
$ws.Range("E5").Value = "HotelLocation.csv"
$ws.Range("F5").Value = "HotelLocation.csv"
$ws.Range("G5").Value = "HotelLocation.csv"

# ---------------------------------------------------------------------------
# Match the formatting used by the rest of the table: plain bordered cells
# for A/B/C/E/F/G/H (same look as row 2 / row 4 non-wrapped cells), and the
# bordered + wrap-text style for the D (formula) column, same as row 4's D
# cell.
# ---------------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)
$ws.Range("E5:H5").PasteSpecial(-4122)

$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Update the view to point at the newly edited row, like Excel does after an
# interactive edit.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("G5").Select()

Write-Output "hotelLocationFileName row inserted"
